# The document has several paragraphs of "boilerplate" course text whose
# contents get cyclically rotated by the diff: each paragraph's <w:t>/
# <w:br/> run content is replaced by text that used to live further down
# the document (and the very last block wraps back around to the
# "Docente" bullet paragraph). None of the paragraphs themselves are
# added/removed/reordered and run-level formatting (bold labels, list
# styles, italics, etc.) is untouched by the diff, so we scope each
# Find/Replace to the specific paragraph's Range (by its fixed index -
# the paragraph count never changes) rather than searching the whole
# document; that keeps every replacement unambiguous even though several
# of the "old" strings appear more than once across the whole doc at
# various points while the script runs.

$d = $word.ActiveDocument
$vt = [char]11   # Word's manual line-break char -> serializes as <w:br/>

function Replace-In-Paragraph($index, $find, $replace) {
    $rng = $d.Paragraphs($index).Range
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 1)
    if (-not $ok) {
        throw "Find.Execute failed in paragraph $index for: $find"
    }
}

$programaText = "Propriedades elétricas: condutividade elétrica em metais puros, ligas metálicas e semicondutores,  e supercondutores; Efeito Hall; Lei de Ohm e dependência com a temperatura." + $vt + "Propriedades magnéticas: susceptibilidade magnética e magnetização c.c. Curvas de histerese de materiais magnéticos macios. Medidas de magnetostricção." + $vt + "Propriedades térmicas dos materiais:  expansão térmica."

$bibliografiaText = "HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000." + $vt + "RAYMOND A. SERWAY, CLEMENT J. MOSES, CURT A. MOYER. Modern Physics 3rd Edition,  Cengage Learning, Inc., 2005." + $vt + "SOLYMAR, L.; WALSH, D. Electrical Properties of Materials, Oxford University Press, 2009." + $vt + "NICOLA A. SPALDIN, Magnetic Materials, Fundamentals and Applications, SECOND EDITION, Cambridge University Press, 2011" + $vt + "ROBERT, P. Electrical and Magnetic Properties of Materials, Artech House, 1998." + $vt + "SPEYER, R. Thermal Analysis of Materials, CRC Press, 1993."

# Paragraph 6: "Objetivos" body <- old "Programa resumido" body
Replace-In-Paragraph 6 "Apresentar as técnicas experimentais de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais." "Estudo das técnicas de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais."

# Paragraph 9: "Docente(s) Responsável(eis)" bullet <- old "Objetivos" body
Replace-In-Paragraph 9 "5840726 - Cristina Bormio Nunes" "Apresentar as técnicas experimentais de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais."

# Paragraph 11: "Programa resumido" body <- old "Programa" body (3 lines, 2 breaks)
Replace-In-Paragraph 11 "Estudo das técnicas de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais." $programaText

# Paragraph 14: "Programa" body <- old "Método" value (collapsing 3 lines into 1)
Replace-In-Paragraph 14 $programaText "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."

# Paragraph 17: "Avaliação" bullet - three separate runs get shifted one slot
# each (Método's value <- old Critério value <- old Norma value <- old
# Bibliografia body). Do these back-to-front: the "Norma" value is
# replaced first, *then* "Critério" is rewritten to what used to be
# Norma's text, *then* "Método" to what used to be Critério's text. Going
# front-to-back would make an earlier rewrite (e.g. Método's new value)
# collide with - and get hit by - a later Find for that same old string.
Replace-In-Paragraph 17 "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação" $bibliografiaText
Replace-In-Paragraph 17 "Média aritmética das notas dos relatórios de cada experimento" "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
Replace-In-Paragraph 17 "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo." "Média aritmética das notas dos relatórios de cada experimento"

# Paragraph 19: "Bibliografia" body <- old "Docente(s)" bullet value (collapsing 6 lines into 1)
Replace-In-Paragraph 19 $bibliografiaText "5840726 - Cristina Bormio Nunes"

Write-Output "done"
